# Add a new "Save" column (H) to the s_vals sheet, matching the header
# style used by the existing columns (B1:G1) and adding a 0 data value
# under it in row 2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the last existing header cell (G1, bold/bordered/
# centered) onto the new header cell H1, then set its text.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H1").Value = "Save"

# New data value for the Save column.
$ws.Range("H2").Value = 0
